$wb = $excel.ActiveWorkbook
$scratch = $wb.Worksheets.Item("Neodymium").Range("ZZ1")

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = 0
$scratch.Formula = "=4.09451772048666E-05"
$ws.Range("C2").Value = $scratch.Value2
$ws.Range("D2").Value = 0.1349277125771417
$ws.Range("E2").Value = 0.8983061480446202
$scratch.Formula = "=3.676266793030401E-12"
$ws.Range("B3").Value = $scratch.Value2
$ws.Range("C3").Value = 0.001952774324661511
$ws.Range("D3").Value = 0.1186635852848033
$ws.Range("E3").Value = 0.7537747045440015
$scratch.Formula = "=5.738750001259246E-14"
$ws.Range("B4").Value = $scratch.Value2
$ws.Range("C4").Value = 0.00176446200604388
$ws.Range("D4").Value = 0.09892343671077571
$ws.Range("E4").Value = 0.6079271585069923
$scratch.Formula = "=3.982554628300847E-08"
$ws.Range("C5").Value = $scratch.Value2
$ws.Range("D5").Value = 0.005361046775579028
$ws.Range("E5").Value = 0.04684562958328455

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030
$scratch.Formula = "=4.827225364041659E-05"
$ws.Range("C2").Value = $scratch.Value2
$ws.Range("D2").Value = 0.1211983383268814
$ws.Range("E2").Value = 1.059056650510759
$ws.Range("C3").Value = 0.002302220284232882
$ws.Range("D3").Value = 0.10658914378472
$ws.Range("E3").Value = 0.8886615276671311
$ws.Range("C4").Value = 0.002080209766060187
$ws.Range("D4").Value = 0.08885762547909347
$ws.Range("E4").Value = 0.7167147877640474
$scratch.Formula = "=4.695226648849531E-08"
$ws.Range("C5").Value = $scratch.Value2
$ws.Range("D5").Value = 0.004815541214495818
$ws.Range("E5").Value = 0.05522858288962464

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030
$scratch.Formula = "=6.274753108837452E-06"
$ws.Range("B2").Value = $scratch.Value2
$ws.Range("C2").Value = 0.005732662478034578
$ws.Range("D2").Value = 0.9345146700513109
$ws.Range("E2").Value = 1.468446836336422
$scratch.Formula = "=4.265998591901206E-05"
$ws.Range("B3").Value = $scratch.Value2
$ws.Range("C3").Value = 0.02071346471615799
$ws.Range("D3").Value = 0.6889245322885504
$ws.Range("E3").Value = 1.05851380799705
$ws.Range("B4").Value = 0.0001264939349260171
$ws.Range("C4").Value = 0.005578321742793286
$ws.Range("D4").Value = 0.580451187597384
$ws.Range("E4").Value = 0.9219790427795368
$scratch.Formula = "=3.974642869529463E-05"
$ws.Range("B5").Value = $scratch.Value2
$ws.Range("C5").Value = 0.01220582660064061
$ws.Range("D5").Value = 0.84504708782342
$ws.Range("E5").Value = 1.211814485320632

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030
$scratch.Formula = "=1.076099402252142E-06"
$ws.Range("B2").Value = $scratch.Value2
$ws.Range("C2").Value = 0.001118055767396142
$ws.Range("D2").Value = 0.5473654256914889
$ws.Range("E2").Value = 1.499733990636409
$scratch.Formula = "=1.148298763755546E-06"
$ws.Range("B3").Value = $scratch.Value2
$ws.Range("C3").Value = 0.00371658900214459
$ws.Range("D3").Value = 0.3007155961130537
$ws.Range("E3").Value = 0.8463127989509267
$scratch.Formula = "=7.353957242698485E-06"
$ws.Range("B4").Value = $scratch.Value2
$ws.Range("C4").Value = 0.001049615036840165
$ws.Range("D4").Value = 0.322155078983428
$ws.Range("E4").Value = 1.057973328863776
$scratch.Formula = "=3.950110178780338E-06"
$ws.Range("B5").Value = $scratch.Value2
$ws.Range("C5").Value = 0.001332061493519333
$ws.Range("D5").Value = 0.5277473432987078
$ws.Range("E5").Value = 1.269697946437359

$scratch.ClearContents()
